$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column C (Forest area (ha)) slightly
$ws.Columns.Item(3).ColumnWidth = 20.83

# Updated data rows 2-22: column A (Number), a new column B (Area ha)
# value is inserted (pushing old B into C), and H/I/J refreshed forecasts.
$ws.Range("A2").Value = 300
$ws.Range("B2").Value = 76657
$ws.Range("C2").Value = 9210
$ws.Range("H2").Value = -35281
$ws.Range("I2").Value = -23822
$ws.Range("J2").Value = 165
$ws.Range("A3").Value = 46
$ws.Range("B3").Value = 2868
$ws.Range("C3").Value = 646
$ws.Range("H3").Value = 17843
$ws.Range("I3").Value = 65324
$ws.Range("J3").Value = 241
$ws.Range("A4").Value = 126
$ws.Range("B4").Value = 6710
$ws.Range("C4").Value = 898
$ws.Range("H4").Value = 12815
$ws.Range("I4").Value = 25901
$ws.Range("J4").Value = 158
$ws.Range("A5").Value = 720
$ws.Range("B5").Value = 119922
$ws.Range("C5").Value = 9564
$ws.Range("H5").Value = 111219
$ws.Range("I5").Value = 201158
$ws.Range("J5").Value = 565
$ws.Range("A6").Value = 959
$ws.Range("B6").Value = 276787
$ws.Range("C6").Value = 34271
$ws.Range("H6").Value = 133064
$ws.Range("I6").Value = 266609
$ws.Range("J6").Value = 622
$ws.Range("A7").Value = 912
$ws.Range("B7").Value = 198551
$ws.Range("C7").Value = 58652
$ws.Range("H7").Value = 279449
$ws.Range("I7").Value = 552482
$ws.Range("J7").Value = 792
$ws.Range("A8").Value = 453
$ws.Range("B8").Value = 59529
$ws.Range("C8").Value = 10614
$ws.Range("H8").Value = 30074
$ws.Range("I8").Value = 70881
$ws.Range("J8").Value = 543
$ws.Range("A9").Value = 336
$ws.Range("B9").Value = 47835
$ws.Range("C9").Value = 26795
$ws.Range("H9").Value = 66802
$ws.Range("I9").Value = 139304
$ws.Range("J9").Value = 490
$ws.Range("A10").Value = 237
$ws.Range("B10").Value = 23477
$ws.Range("C10").Value = 5100
$ws.Range("H10").Value = 18931
$ws.Range("I10").Value = 51405
$ws.Range("J10").Value = 278
$ws.Range("A11").Value = 529
$ws.Range("B11").Value = 69644
$ws.Range("C11").Value = 8732
$ws.Range("H11").Value = 13917
$ws.Range("I11").Value = 39339
$ws.Range("J11").Value = 256
$ws.Range("A12").Value = 415
$ws.Range("B12").Value = 107979
$ws.Range("C12").Value = 61888
$ws.Range("H12").Value = -9778
$ws.Range("I12").Value = -4962
$ws.Range("J12").Value = 581
$ws.Range("A13").Value = 825
$ws.Range("B13").Value = 188100
$ws.Range("C13").Value = 93465
$ws.Range("H13").Value = 83766
$ws.Range("I13").Value = 204556
$ws.Range("J13").Value = 792
$ws.Range("A14").Value = 1235
$ws.Range("B14").Value = 1317798
$ws.Range("C14").Value = 786198
$ws.Range("H14").Value = 256124
$ws.Range("I14").Value = 465578
$ws.Range("J14").Value = 738
$ws.Range("A15").Value = 517
$ws.Range("B15").Value = 200479
$ws.Range("C15").Value = 134647
$ws.Range("H15").Value = 108355
$ws.Range("I15").Value = 214554
$ws.Range("J15").Value = 560
$ws.Range("A16").Value = 94
$ws.Range("B16").Value = 6074
$ws.Range("C16").Value = 1775
$ws.Range("H16").Value = -39082
$ws.Range("I16").Value = -70147
$ws.Range("J16").Value = 478
$ws.Range("A17").Value = 122
$ws.Range("B17").Value = 10418
$ws.Range("C17").Value = 2338
$ws.Range("H17").Value = 52791
$ws.Range("I17").Value = 113328
$ws.Range("J17").Value = 423
$ws.Range("A18").Value = 244
$ws.Range("B18").Value = 28917
$ws.Range("C18").Value = 18388
$ws.Range("H18").Value = 226122
$ws.Range("I18").Value = 360637
$ws.Range("J18").Value = 747
$ws.Range("A19").Value = 383
$ws.Range("B19").Value = 99483
$ws.Range("C19").Value = 81724
$ws.Range("H19").Value = 28372
$ws.Range("I19").Value = 67221
$ws.Range("J19").Value = 404
$ws.Range("A20").Value = 558
$ws.Range("B20").Value = 24722
$ws.Range("C20").Value = 12674
$ws.Range("H20").Value = 9783
$ws.Range("I20").Value = 30651
$ws.Range("J20").Value = 299
$ws.Range("A21").Value = 543
$ws.Range("B21").Value = 31199
$ws.Range("C21").Value = 9198
$ws.Range("H21").Value = 18929
$ws.Range("I21").Value = 47601
$ws.Range("J21").Value = 415
$ws.Range("A22").Value = 633
$ws.Range("B22").Value = 227985
$ws.Range("C22").Value = 185155
$ws.Range("H22").Value = 167720
$ws.Range("I22").Value = 307538
$ws.Range("J22").Value = 639
